# Destajo Sem 40 sin tablas y totales c/formato
#
# Appends the Semana (week) 40 destajo entries as new rows to the "Tabla3"
# Excel table on Hoja1. Each hashtable below is one row of data keyed by
# column letter (A..V), matching the table's column layout:
#   A=SEMANA  C=FECHA  D=INSTALADOR  E=CODIGO  F=PROYECTO  G=ML  H=$ X ML
#   I=TOTAL ML  J=PZAS  K=$ X PZA  L=TOTAL PZAS  M=DIA  N=$ X DIA
#   O=TOTAL LOTE  P=SUBTOTAL  Q=BONO PUNTUALIDAD  S=MONTO  U=TOTAL DESTAJO

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

$data = @(
    [ordered]@{'A'=40; 'C'=45568; 'D'='Elías Cadpevila Figueroa'; 'E'='CTVV 2502'; 'F'='CTVV 2502-CABO SUR PROTOTIPO MALDOVA'; 'I'=0; 'L'=0; 'M'=1; 'N'=750; 'O'=750; 'P'=750; 'S'=750; 'U'=750},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Marco Antonio Barrera Bautista'; 'E'='CTVV 2502'; 'F'='CTVV 2502-CABO SUR PROTOTIPO MALDOVA'; 'I'=0; 'L'=0; 'M'=0.5; 'N'=700; 'O'=350; 'P'=350; 'S'=350; 'U'=350},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Marco Antonio Barrera Bautista'; 'E'='CTVV 2431'; 'F'='CTVV 2431 - Hercom (3era etapa) - Hercom'; 'G'=0.85; 'H'=700; 'I'=595; 'J'=1; 'K'=365; 'L'=365; 'M'=1; 'N'=750; 'O'=750; 'P'=1710; 'S'=1710; 'U'=1710},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Marco Antonio Barrera Bautista'; 'E'='CTVV 2431'; 'F'='CTVV 2431 - Hercom (3era etapa) - Hercom'; 'G'=2.59; 'H'=350; 'I'=906.5; 'J'=1; 'K'=700; 'L'=700; 'O'=0; 'P'=1606.5; 'S'=1606.5; 'U'=1606.5},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Marco Antonio Barrera Bautista'; 'E'='CTVV 2431'; 'F'='CTVV 2431 - Hercom (3era etapa) - Hercom'; 'G'=0.71; 'H'=700; 'I'=497; 'L'=0; 'O'=0; 'P'=497; 'Q'=500; 'S'=997; 'U'=997},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Jorge Armando Calderón'; 'E'='CTVV 2529'; 'F'='CTVV 2529 - VIALBA PROTOTIPO T01-18 CASAS - GRUPO HERSO'; 'G'=1.0900000000000001; 'H'=320; 'I'=348.8; 'J'=2; 'K'=300; 'L'=600; 'O'=0; 'P'=948.8; 'S'=948.8; 'U'=948.8},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Jorge Armando Calderón'; 'E'='CTVV 2529'; 'F'='CTVV 2529 - VIALBA PROTOTIPO T01-18 CASAS - GRUPO HERSO'; 'G'=0.81; 'H'=700; 'I'=567; 'J'=1; 'K'=150; 'L'=150; 'O'=0; 'P'=717; 'S'=717; 'U'=717},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Jorge Armando Calderón'; 'E'='CTVV 2529'; 'F'='CTVV 2529 - VIALBA PROTOTIPO T01-18 CASAS - GRUPO HERSO'; 'G'=1.0900000000000001; 'H'=320; 'I'=348.8; 'L'=0; 'O'=0; 'P'=348.8; 'S'=348.8; 'U'=348.8},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Jorge Armando Calderón'; 'E'='CTVV 2529'; 'F'='CTVV 2529 - VIALBA PROTOTIPO T01-18 CASAS - GRUPO HERSO'; 'G'=1.9; 'H'=700; 'I'=1330; 'L'=0; 'O'=0; 'P'=1330; 'S'=1330; 'U'=1330},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Jorge Armando Calderón'; 'E'='CTVV 2529'; 'F'='CTVV 2529 - VIALBA PROTOTIPO T01-18 CASAS - GRUPO HERSO'; 'G'=1.0900000000000001; 'H'=320; 'I'=348.8; 'L'=0; 'O'=0; 'P'=348.8; 'S'=348.8; 'U'=348.8},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Jorge Armando Calderón'; 'E'='CTVV 2529'; 'F'='CTVV 2529 - VIALBA PROTOTIPO T01-18 CASAS - GRUPO HERSO'; 'G'=1.9; 'H'=700; 'I'=1330; 'L'=0; 'O'=0; 'P'=1330; 'S'=1330; 'U'=1330},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Julio Cesar Domínguez'; 'E'='CTVV 2651'; 'F'='CTVV 2651 - Veronesa Prototipo A - Constructora GyH'; 'G'=2.27; 'H'=320; 'I'=726.4; 'J'=1; 'K'=500; 'L'=500; 'O'=0; 'P'=1226.4000000000001; 'S'=1226.4000000000001; 'U'=1226.4000000000001},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Jesus Soto Días'; 'E'='CTVV 3142'; 'F'='CTVV 3142 - Propuesta Verona Zamora (10 casas) - Meda'; 'G'=2.16; 'H'=700; 'I'=1512; 'L'=0; 'O'=0; 'P'=1512; 'S'=1512; 'U'=1512},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Jesus Soto Días'; 'E'='CTVV 3142'; 'F'='CTVV 3142 - Propuesta Verona Zamora (10 casas) - Meda'; 'G'=2.8; 'H'=700; 'I'=1959.9999999999998; 'L'=0; 'O'=0; 'P'=1959.9999999999998; 'S'=1959.9999999999998; 'U'=1959.9999999999998},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Jesus Soto Días'; 'E'='CTVV 3142'; 'F'='CTVV 3142 - Propuesta Verona Zamora (10 casas) - Meda'; 'G'=2.16; 'H'=700; 'I'=1512; 'L'=0; 'O'=0; 'P'=1512; 'S'=1512; 'U'=1512},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Carlos Alfredo Reyes Hernández'; 'E'='CTVV 2776'; 'F'='CTVV 2776 - Casa JML Primera Etapa - Miguel Maldonado'; 'I'=0; 'L'=0; 'M'=5; 'N'=1083; 'O'=5415; 'P'=5415; 'S'=5415; 'U'=5415},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Luis Alberto Ponce'; 'E'='CTVV 3144'; 'F'='CTVV 3144 - Torre Aria Prototipo A (2 departamentos) - Injesa'; 'G'=2.6; 'H'=200; 'I'=520; 'J'=5; 'K'=700; 'L'=3500; 'O'=0; 'P'=4020; 'S'=4020; 'U'=4020},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Luis Alberto Ponce'; 'E'='CTVV 3144'; 'F'='CTVV 3144 - Torre Aria Prototipo A (2 departamentos) - Injesa'; 'G'=2.4500000000000002; 'H'=200; 'I'=490.00000000000006; 'L'=0; 'O'=0; 'P'=490.00000000000006; 'S'=490.00000000000006; 'U'=490.00000000000006},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Agustín Árias Venegas'; 'E'='CTVV 3078'; 'F'='CTVV 3078 - Piso 5 Centro Capital Oficina 02 LN - Canaco'; 'G'=0.41; 'H'=900; 'I'=369; 'L'=0; 'M'=0.5; 'N'=1083.33; 'O'=541.66499999999996; 'P'=910.66499999999996; 'S'=910.66499999999996; 'U'=910.66499999999996},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Agustín Árias Venegas'; 'E'='CTVV 3078'; 'F'='CTVV 3078 - Piso 5 Centro Capital Oficina 02 LN - Canaco'; 'G'=0.64; 'H'=300; 'I'=192; 'L'=0; 'O'=0; 'P'=192; 'S'=192; 'U'=192},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Agustín Árias Venegas'; 'E'='CTVV 3078'; 'F'='CTVV 3078 - Piso 5 Centro Capital Oficina 02 LN - Canaco'; 'G'=2; 'H'=300; 'I'=600; 'L'=0; 'O'=0; 'P'=600; 'S'=600; 'U'=600},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Agustín Árias Venegas'; 'E'='CTVV 3078'; 'F'='CTVV 3078 - Piso 5 Centro Capital Oficina 02 LN - Canaco'; 'G'=2.76; 'H'=300; 'I'=827.99999999999989; 'L'=0; 'O'=0; 'P'=827.99999999999989; 'S'=827.99999999999989; 'U'=827.99999999999989},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Agustín Árias Venegas'; 'E'='CTVV 3078'; 'F'='CTVV 3078 - Piso 5 Centro Capital Oficina 02 LN - Canaco'; 'G'=2.76; 'H'=100; 'I'=276; 'L'=0; 'O'=0; 'P'=276; 'S'=276; 'U'=276},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Luis Alberto Andrade Martínez'; 'E'='CTVV 3078'; 'F'='CTVV 3078 - Piso 5 Centro Capital Oficina 02 LN - Canaco'; 'G'=3.08; 'H'=400; 'I'=1232; 'L'=0; 'O'=0; 'P'=1232; 'Q'=500; 'S'=1732; 'U'=1732},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Luis Alberto Andrade Martínez'; 'E'='CTVV 2288'; 'F'='CTVV 2288 - VISTA VERDE PROTOTIPO SAUCE'; 'I'=0; 'J'=1; 'K'=250; 'L'=250; 'O'=0; 'P'=250; 'S'=250; 'U'=250},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Luis Alberto Andrade Martínez'; 'E'='CTVV 2288'; 'F'='CTVV 2288 - VISTA VERDE PROTOTIPO SAUCE'; 'I'=0; 'J'=1; 'K'=750; 'L'=750; 'O'=0; 'P'=750; 'S'=750; 'U'=750},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Luis Alberto Andrade Martínez'; 'E'='CTVV 2288'; 'F'='CTVV 2288 - VISTA VERDE PROTOTIPO SAUCE'; 'I'=0; 'J'=1; 'K'=1800; 'L'=1800; 'O'=0; 'P'=1800; 'S'=1800; 'U'=1800},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Luis Alberto Andrade Martínez'; 'E'='CTVV 2288'; 'F'='CTVV 2288 - VISTA VERDE PROTOTIPO SAUCE'; 'G'=10.5; 'H'=250; 'I'=2625; 'L'=0; 'O'=0; 'P'=2625; 'S'=2625; 'U'=2625},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Luis Alberto Andrade Martínez'; 'E'='CTVV 2993'; 'F'='CTVV 2993 - Prototipo Bilbao Torres Castillejo - Hogares Deesa'; 'G'=3.8650000000000002; 'H'=700; 'I'=2705.5; 'J'=1; 'K'=400; 'L'=400; 'O'=0; 'P'=3105.5; 'S'=3105.5; 'U'=3105.5},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Luis Alberto Andrade Martínez'; 'E'='CTVV 2993'; 'F'='CTVV 2993 - Prototipo Bilbao Torres Castillejo - Hogares Deesa'; 'G'=3.8650000000000002; 'H'=700; 'I'=2705.5; 'J'=1; 'K'=400; 'L'=400; 'O'=0; 'P'=3105.5; 'S'=3105.5; 'U'=3105.5},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Luis Alberto Andrade Martínez'; 'E'='CTVV 2993'; 'F'='CTVV 2993 - Prototipo Bilbao Torres Castillejo - Hogares Deesa'; 'G'=1.8049999999999999; 'H'=700; 'I'=1263.5; 'J'=1; 'K'=200; 'L'=200; 'O'=0; 'P'=1463.5; 'S'=1463.5; 'U'=1463.5},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Luis Alberto Andrade Martínez'; 'E'='CTVV 2993'; 'F'='CTVV 2993 - Prototipo Bilbao Torres Castillejo - Hogares Deesa'; 'G'=2.3199999999999998; 'H'=700; 'I'=1624; 'J'=1; 'K'=200; 'L'=200; 'O'=0; 'P'=1824; 'S'=1824; 'U'=1824},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Abel Guzmán García'; 'E'='CTVV 2993'; 'F'='CTVV 2993 - Prototipo Bilbao Torres Castillejo - Hogares Deesa'; 'G'=3.09; 'H'=700; 'I'=2163; 'L'=0; 'O'=0; 'P'=2163; 'S'=2163; 'U'=2163},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Jorge Zamudio Calderón'; 'E'='CTVV 3125'; 'F'='CTVV 3125 - Remodelación Puerta Sur - Verónica Erón'; 'I'=0; 'J'=1; 'K'=750; 'L'=750; 'O'=0; 'P'=750; 'S'=750; 'U'=750},
    [ordered]@{'A'=40; 'C'=45568; 'D'='Juan Luis Ramírez'; 'E'='CTVV 2109'; 'F'='CTVV 2109 - CASA SO MODISA - MODISA'; 'I'=0; 'L'=0; 'M'=1; 'N'=5000; 'O'=5000; 'P'=5000; 'S'=5000; 'U'=5000}
)

# Column letter -> 1-based worksheet column index (A..V)
$colIndex = @{
    'A'=1;  'B'=2;  'C'=3;  'D'=4;  'E'=5;  'F'=6;  'G'=7;  'H'=8;
    'I'=9;  'J'=10; 'K'=11; 'L'=12; 'M'=13; 'N'=14; 'O'=15; 'P'=16;
    'Q'=17; 'R'=18; 'S'=19; 'T'=20; 'U'=21; 'V'=22
}

# Row number (1-based, sheet-wide) of the first row we are about to add.
$firstNewRow = $lo.Range.Row + $lo.Range.Rows.Count

# ListRows.Add() appends one blank row at the bottom of the table and
# automatically grows the table ref, the autoFilter and the sheet
# dimension to match - same as typing into the row right below an
# Excel Table in the UI.
foreach ($rowData in $data) {
    $lo.ListRows.Add() | Out-Null
}

# Now fill in the values for every new row.
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $firstNewRow + $i
    $rowData = $data[$i]
    foreach ($col in $rowData.Keys) {
        $ws.Cells.Item($rowNum, $colIndex[$col]).Value = $rowData[$col]
    }
}

# Match the author's final on-screen selection: the newly appended block.
$lastRow = $firstNewRow + $data.Count - 1
$ws.Activate()
$ws.Range("A$($firstNewRow):V$($lastRow)").Select()
